$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.92
$wsSummary.Range("B4").Value = -0.08
$wsSummary.Range("B5").Value = -0.27
$wsSummary.Range("B6").Value = 6
$wsSummary.Range("B7").Value = 3
$wsSummary.Range("B9").Value = 50

# --- Strategy Status sheet ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.92
$wsStatus.Range("D4").Value = 6
$wsStatus.Range("E4").Value = -0.08
$wsStatus.Range("F4").Value = -0.08
$wsStatus.Range("G4").Value = 50

# --- New trade row data (Trade #6) to append to both "All Trades" and "MarketMaking" sheets ---
function Add-TradeRow($ws) {
    $ws.Cells.Item(7, 1).Value = 6
    # The Date column already holds the identical "2026-02-17" text in every
    # existing row, so copy it down rather than re-assigning the value -
    # assigning a date-shaped string directly makes Excel auto-convert it to
    # a date serial number instead of keeping it as plain text.
    $ws.Cells.Item(2, 2).Copy($ws.Cells.Item(7, 2))
    # The Time column value is new; assigning it via Formula (rather than
    # Value) avoids Excel's date/time auto-conversion and keeps it text.
    $ws.Cells.Item(7, 3).Formula = "08:07:57"
    $ws.Cells.Item(7, 4).Value = "MarketMaking"
    $ws.Cells.Item(7, 5).Value = "DOWN"
    $ws.Cells.Item(7, 6).Value = 0.66
    $ws.Cells.Item(7, 7).Value = 0.7
    $ws.Cells.Item(7, 8).Value = "CLOSED"
    $ws.Cells.Item(7, 9).Value = 6.0606
    $ws.Cells.Item(7, 10).Value = 0.04
    $ws.Cells.Item(7, 11).Value = 99.92
    $ws.Cells.Item(7, 12).Value = 0
    $ws.Cells.Item(7, 13).Value = 0
    $ws.Cells.Item(7, 14).Value = 0.6
    $ws.Cells.Item(7, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(7, 16).Value = "early_exit"
    $ws.Cells.Item(7, 17).Value = 0.13
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
